$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.145.11'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -3.26%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.861.84'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -4.30%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9996'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.26%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '233.60'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -3.71%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9995'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.24%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4648'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -3.66%  '

$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -3.37%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06551'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -3.73%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.86'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.13%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07809'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.62%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '96.58'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -7.58%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.864.25'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -4.24%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.132'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -3.40%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6673'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -3.39%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '281.56'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -5.18%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '30.181.88'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -3.13%  '

$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.12%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.522'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -1.21%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.62'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -3.03%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.103.21'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -5.05%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.000007240'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -4.99%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9995'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.33%  '

$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -5.17%  '

$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -2.68%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '166.31'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -1.67%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.89'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -4.88%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.911'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -10.88%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.340'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -4.00%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.09567'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -5.87%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.406'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -4.85%  '

$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -4.48%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.111'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -5.82%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04656'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -3.91%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7011'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -5.62%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.093'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -3.64%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.699'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -1.71%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01852'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -5.57%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.276'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -5.08%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.517'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -4.93%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '73.58'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -5.04%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8555'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -2.23%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.919'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -5.64%  '

$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value = 'PaxDollar'
$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9989'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.31%  '

$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'TheSandbox'
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4155'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -4.98%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '103.59'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -2.78%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '995.61'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -2.55%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.193'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -5.10%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.213'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.23%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '34.15'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -2.99%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.1139'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -6.28%  '
